$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 123.666664
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = $null
$ws.Range("H4").Value = 812.25
$ws.Range("I4").Value = 812.25
$ws.Range("K4").Value = 812.25
$ws.Range("M4").Value = -698.25
$ws.Range("H74").Value = 5517.8887
$ws.Range("J74").Value = 5250
$ws.Range("L74").Value = 5250
$ws.Range("N74").Value = -7122
$ws.Range("H77").Value = 5517.8887
$ws.Range("J77").Value = 5250
$ws.Range("L77").Value = 26250
$ws.Range("N77").Value = -35610
$ws.Range("H111").Value = 2102.3333
$ws.Range("I111").Value = 2102.3333
$ws.Range("K111").Value = 6306.999899999999
$ws.Range("M111").Value = -3239.999899999999
$ws.Range("H132").Value = 2441.4167
$ws.Range("I132").Value = 1929.8
$ws.Range("K132").Value = 5789.4
$ws.Range("M132").Value = -3259.4
$ws.Range("H141").Value = 4368
$ws.Range("I141").Value = 4595.4
$ws.Range("K141").Value = 13786.2
$ws.Range("M141").Value = -8606.199999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4728.914
$ws.Range("J2").Value = 8170.7144
$ws.Range("L2").Value = 8170.7144
$ws.Range("N2").Value = -8396.714400000001
$ws.Range("H4").Value = 2166.6667
$ws.Range("I4").Value = 2250
$ws.Range("K4").Value = 2250
$ws.Range("M4").Value = -2134
$ws.Range("H33").Value = 18713
$ws.Range("I33").Value = 16617.334
$ws.Range("J33").Value = 25000
$ws.Range("K33").Value = 16617.334
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = -16288.334
$ws.Range("N33").Value = -25658
$ws.Range("H36").Value = 26833.334
$ws.Range("J36").Value = 25000
$ws.Range("L36").Value = 25000
$ws.Range("N36").Value = -25692
$ws.Range("H61").Value = 4965.5107
$ws.Range("I61").Value = 3320.4866
$ws.Range("K61").Value = 3320.4866
$ws.Range("M61").Value = -3108.4866
$ws.Range("H74").Value = 4256.2
$ws.Range("I74").Value = 3817.923
$ws.Range("J74").Value = 4855.9473
$ws.Range("K74").Value = 3817.923
$ws.Range("L74").Value = 4855.9473
$ws.Range("M74").Value = -2943.923
$ws.Range("N74").Value = -6603.9473
$ws.Range("H77").Value = 4256.2
$ws.Range("I77").Value = 3817.923
$ws.Range("J77").Value = 4855.9473
$ws.Range("K77").Value = 19089.615
$ws.Range("L77").Value = 24279.7365
$ws.Range("M77").Value = -14721.615
$ws.Range("N77").Value = -33015.7365
$ws.Range("H116").Value = 4728.914
$ws.Range("J116").Value = 8170.7144
$ws.Range("L116").Value = 8170.7144
$ws.Range("N116").Value = -12758.7144
$ws.Range("H132").Value = 3742.5527
$ws.Range("I132").Value = 3623.1
$ws.Range("K132").Value = 10869.3
$ws.Range("M132").Value = -8339.299999999999
$ws.Range("H136").Value = 4965.5107
$ws.Range("I136").Value = 3320.4866
$ws.Range("K136").Value = 9961.459800000001
$ws.Range("M136").Value = -7411.459800000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4728.914
$ws.Range("J3").Value = 8170.7144
$ws.Range("L3").Value = 8170.7144
$ws.Range("N3").Value = -8398.714400000001
$ws.Range("H86").Value = 325504.75
$ws.Range("I86").Value = 529169.4399999999
$ws.Range("K86").Value = 529169.4399999999
$ws.Range("M86").Value = -528046.4399999999
$ws.Range("H89").Value = 325504.75
$ws.Range("I89").Value = 529169.4399999999
$ws.Range("K89").Value = 2645847.2
$ws.Range("M89").Value = -2640231.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2860821.5
$ws.Range("I4").Value = 168.66667
$ws.Range("J4").Value = 3640999.5
$ws.Range("K4").Value = 168.66667
$ws.Range("L4").Value = 3640999.5
$ws.Range("M4").Value = -56.66667000000001
$ws.Range("N4").Value = -3641223.5
$ws.Range("H58").Value = 4743.3706
$ws.Range("I58").Value = 1717.2858
$ws.Range("J58").Value = 8002.231
$ws.Range("K58").Value = 1717.2858
$ws.Range("L58").Value = 8002.231
$ws.Range("M58").Value = -1514.2858
$ws.Range("N58").Value = -8408.231
$ws.Range("H62").Value = 98043.09
$ws.Range("J62").Value = 8559.5
$ws.Range("L62").Value = 8559.5
$ws.Range("N62").Value = -9807.5
$ws.Range("H65").Value = 98043.09
$ws.Range("J65").Value = 8559.5
$ws.Range("L65").Value = 42797.5
$ws.Range("N65").Value = -49037.5
$ws.Range("H122").Value = 1785.75
$ws.Range("I122").Value = 1727
$ws.Range("K122").Value = 5181
$ws.Range("M122").Value = -2731
$ws.Range("H132").Value = 29833.791
$ws.Range("I132").Value = 21485.107
$ws.Range("J132").Value = 41521.95
$ws.Range("K132").Value = 64455.321
$ws.Range("L132").Value = 124565.85
$ws.Range("M132").Value = -61925.321
$ws.Range("N132").Value = -129625.85
$ws.Range("H134").Value = 6494.7393
$ws.Range("I134").Value = 6012.154
$ws.Range("J134").Value = 7122.1
$ws.Range("K134").Value = 18036.462
$ws.Range("L134").Value = 21366.3
$ws.Range("M134").Value = -15501.462
$ws.Range("N134").Value = -26436.3
$ws.Range("H136").Value = 4743.3706
$ws.Range("I136").Value = 1717.2858
$ws.Range("J136").Value = 8002.231
$ws.Range("K136").Value = 5151.857400000001
$ws.Range("L136").Value = 24006.693
$ws.Range("M136").Value = -2601.857400000001
$ws.Range("N136").Value = -29106.693

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 250
$ws.Range("J31").Value = 250
$ws.Range("L31").Value = 750
$ws.Range("N31").Value = -1326
$ws.Range("H75").Value = 2482.2856
$ws.Range("J75").Value = 2435.2
$ws.Range("L75").Value = 7305.599999999999
$ws.Range("N75").Value = -9301.599999999999
$ws.Range("H78").Value = 2482.2856
$ws.Range("J78").Value = 2435.2
$ws.Range("L78").Value = 21916.8
$ws.Range("N78").Value = -31900.8
$ws.Range("H117").Value = 714.625
$ws.Range("I117").Value = 819.5
$ws.Range("J117").Value = 400
$ws.Range("K117").Value = 2458.5
$ws.Range("L117").Value = 1200
$ws.Range("M117").Value = 983.5
$ws.Range("N117").Value = -8084

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 22647.5
$ws.Range("I5").Value = 29450
$ws.Range("K5").Value = 29450
$ws.Range("M5").Value = -29338
$ws.Range("H122").Value = 2610.7646
$ws.Range("I122").Value = 2323.2593
$ws.Range("J122").Value = 3719.7144
$ws.Range("K122").Value = 6969.777900000001
$ws.Range("L122").Value = 11159.1432
$ws.Range("M122").Value = -4519.777900000001
$ws.Range("N122").Value = -16059.1432
$ws.Range("H132").Value = 9431.643
$ws.Range("I132").Value = 10204.739
$ws.Range("J132").Value = 5875.4
$ws.Range("K132").Value = 30614.217
$ws.Range("L132").Value = 17626.2
$ws.Range("M132").Value = -28084.217
$ws.Range("N132").Value = -22686.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4772.222
$ws.Range("I136").Value = 4535.4585
$ws.Range("K136").Value = 13606.3755
$ws.Range("M136").Value = -11056.3755

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = $null
$ws.Range("H132").Value = 1752804.6
$ws.Range("I132").Value = 2654613.5
$ws.Range("K132").Value = 7963840.5
$ws.Range("M132").Value = -7961310.5
$ws.Range("H136").Value = 2066.3333
$ws.Range("I136").Value = 2066.3333
$ws.Range("K136").Value = 6198.999899999999
$ws.Range("M136").Value = -3648.999899999999
